# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Net effect on Hoja1: the two worker blocks (rows 16-17 and 18-19) swap
# so "LUIS ALBERTO MIRANDA MARTINEZ" (1002274459) comes first, followed by
# "ALBERTO ANTONIO VALDEZ VELASQUEZ" (1044927350); values/periods follow
# along with each worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: LUIS ALBERTO MIRANDA MARTINEZ, periodo 2405, valor 73333
$ws.Range("C16").Value = "1002274459"
$ws.Range("D16").Value = "LUIS ALBERTO MIRANDA MARTINEZ"
$ws.Range("E16").Value = "2405"
$ws.Range("F16").Value = 73333

# Row 17: LUIS ALBERTO MIRANDA MARTINEZ, periodo 2404, valor 48889
$ws.Range("C17").Value = "1002274459"
$ws.Range("D17").Value = "LUIS ALBERTO MIRANDA MARTINEZ"
$ws.Range("E17").Value = "2404"
$ws.Range("F17").Value = 48889

# Row 18: ALBERTO ANTONIO VALDEZ VELASQUEZ, periodo 2405, valor 73333
$ws.Range("C18").Value = "1044927350"
$ws.Range("D18").Value = "ALBERTO ANTONIO VALDEZ VELASQUEZ"
$ws.Range("E18").Value = "2405"
$ws.Range("F18").Value = 73333

# Row 19: ALBERTO ANTONIO VALDEZ VELASQUEZ, periodo 2404, valor 48889
$ws.Range("C19").Value = "1044927350"
$ws.Range("D19").Value = "ALBERTO ANTONIO VALDEZ VELASQUEZ"
$ws.Range("E19").Value = "2404"
$ws.Range("F19").Value = 48889
